# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the per-language handback status sheets.

$wb = $excel.ActiveWorkbook

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-19 04:20:39"
$wsZh.Range("E5").Value = "2016-03-19 04:20:39"
$wsZh.Range("H2").Value = "2016-03-19 04:20:58"
$wsZh.Range("H5").Value = "2016-03-19 04:20:58"

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-19 04:20:42"
$wsDe.Range("E5").Value = "2016-03-19 04:20:42"
$wsDe.Range("H2").Value = "2016-03-19 04:21:04"
$wsDe.Range("H5").Value = "2016-03-19 04:21:04"
